# Daily refresh of the Kahraman XGB forecast table.
# Drops the oldest day of hourly predictions and appends the newest
# day (Dragosel) that has just completed a forecast run, shifting the
# whole rolling 169-hour window forward by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is: ExcelDateSerial, Interval (hour-of-day), Prediction, Lookup label
$forecast = @(
    @(45890, 17, 0.107, "21.08.202517"),
    @(45890, 18, 0.107, "21.08.202518"),
    @(45890, 19, 0.791, "21.08.202519"),
    @(45890, 20, 0.296, "21.08.202520"),
    @(45890, 21, 0.016, "21.08.202521"),
    @(45890, 22, 0.031, "21.08.202522"),
    @(45890, 23, 0.026, "21.08.202523"),
    @(45890, 24, 0, "21.08.202524"),
    @(45891, 1, 0, "22.08.20251"),
    @(45891, 2, 0, "22.08.20252"),
    @(45891, 3, 0.013, "22.08.20253"),
    @(45891, 4, 0.015, "22.08.20254"),
    @(45891, 5, 0.015, "22.08.20255"),
    @(45891, 6, 0, "22.08.20256"),
    @(45891, 7, 0, "22.08.20257"),
    @(45891, 8, 0.174, "22.08.20258"),
    @(45891, 9, 0.756, "22.08.20259"),
    @(45891, 10, 1.512, "22.08.202510"),
    @(45891, 11, 2.601, "22.08.202511"),
    @(45891, 12, 3.127, "22.08.202512"),
    @(45891, 13, 3.34, "22.08.202513"),
    @(45891, 14, 3.401, "22.08.202514"),
    @(45891, 15, 2.848, "22.08.202515"),
    @(45891, 16, 2.566, "22.08.202516"),
    @(45891, 17, 2.005, "22.08.202517"),
    @(45891, 18, 1.612, "22.08.202518"),
    @(45891, 19, 0.722, "22.08.202519"),
    @(45891, 20, 0.171, "22.08.202520"),
    @(45891, 21, 0, "22.08.202521"),
    @(45891, 22, 0, "22.08.202522"),
    @(45891, 23, 0, "22.08.202523"),
    @(45891, 24, 0, "22.08.202524"),
    @(45892, 1, 0, "23.08.20251"),
    @(45892, 2, 0, "23.08.20252"),
    @(45892, 3, 0, "23.08.20253"),
    @(45892, 4, 0, "23.08.20254"),
    @(45892, 5, 0, "23.08.20255"),
    @(45892, 6, 0, "23.08.20256"),
    @(45892, 7, 0, "23.08.20257"),
    @(45892, 8, 0.353, "23.08.20258"),
    @(45892, 9, 1.408, "23.08.20259"),
    @(45892, 10, 2.488, "23.08.202510"),
    @(45892, 11, 3.387, "23.08.202511"),
    @(45892, 12, 3.747, "23.08.202512"),
    @(45892, 13, 3.848, "23.08.202513"),
    @(45892, 14, 3.59, "23.08.202514"),
    @(45892, 15, 3.459, "23.08.202515"),
    @(45892, 16, 3.315, "23.08.202516"),
    @(45892, 17, 2.694, "23.08.202517"),
    @(45892, 18, 1.804, "23.08.202518"),
    @(45892, 19, 1.011, "23.08.202519"),
    @(45892, 20, 0.302, "23.08.202520"),
    @(45892, 21, 0.028, "23.08.202521"),
    @(45892, 22, 0, "23.08.202522"),
    @(45892, 23, 0.013, "23.08.202523"),
    @(45892, 24, 0.012, "23.08.202524"),
    @(45893, 1, 0.011, "24.08.20251"),
    @(45893, 2, 0.011, "24.08.20252"),
    @(45893, 3, 0, "24.08.20253"),
    @(45893, 4, 0, "24.08.20254"),
    @(45893, 5, 0, "24.08.20255"),
    @(45893, 6, 0, "24.08.20256"),
    @(45893, 7, 0, "24.08.20257"),
    @(45893, 8, 0.276, "24.08.20258"),
    @(45893, 9, 1.074, "24.08.20259"),
    @(45893, 10, 1.877, "24.08.202510"),
    @(45893, 11, 2.471, "24.08.202511"),
    @(45893, 12, 2.711, "24.08.202512"),
    @(45893, 13, 2.788, "24.08.202513"),
    @(45893, 14, 2.989, "24.08.202514"),
    @(45893, 15, 3.339, "24.08.202515"),
    @(45893, 16, 3.066, "24.08.202516"),
    @(45893, 17, 2.762, "24.08.202517"),
    @(45893, 18, 1.887, "24.08.202518"),
    @(45893, 19, 1.138, "24.08.202519"),
    @(45893, 20, 0.356, "24.08.202520"),
    @(45893, 21, 0.045, "24.08.202521"),
    @(45893, 22, 0.06900000000000001, "24.08.202522"),
    @(45893, 23, 0.112, "24.08.202523"),
    @(45893, 24, 0.32, "24.08.202524"),
    @(45894, 1, 0.094, "25.08.20251"),
    @(45894, 2, 0.106, "25.08.20252"),
    @(45894, 3, 0.106, "25.08.20253"),
    @(45894, 4, 0.106, "25.08.20254"),
    @(45894, 5, 0.108, "25.08.20255"),
    @(45894, 6, 0.109, "25.08.20256"),
    @(45894, 7, 0.111, "25.08.20257"),
    @(45894, 8, 0.435, "25.08.20258"),
    @(45894, 9, 1.973, "25.08.20259"),
    @(45894, 10, 2.776, "25.08.202510"),
    @(45894, 11, 3.559, "25.08.202511"),
    @(45894, 12, 4.181, "25.08.202512"),
    @(45894, 13, 3.946, "25.08.202513"),
    @(45894, 14, 3.907, "25.08.202514"),
    @(45894, 15, 3.96, "25.08.202515"),
    @(45894, 16, 4.006, "25.08.202516"),
    @(45894, 17, 3.507, "25.08.202517"),
    @(45894, 18, 2.617, "25.08.202518"),
    @(45894, 19, 1.419, "25.08.202519"),
    @(45894, 20, 0.33, "25.08.202520"),
    @(45894, 21, 0.065, "25.08.202521"),
    @(45894, 22, 0.101, "25.08.202522"),
    @(45894, 23, 0.106, "25.08.202523"),
    @(45894, 24, 0.101, "25.08.202524"),
    @(45895, 1, 0.103, "26.08.20251"),
    @(45895, 2, 0.09, "26.08.20252"),
    @(45895, 3, 0.08799999999999999, "26.08.20253"),
    @(45895, 4, 0.08, "26.08.20254"),
    @(45895, 5, 0.08799999999999999, "26.08.20255"),
    @(45895, 6, 0.07099999999999999, "26.08.20256"),
    @(45895, 7, 0.063, "26.08.20257"),
    @(45895, 8, 0.392, "26.08.20258"),
    @(45895, 9, 1.535, "26.08.20259"),
    @(45895, 10, 2.586, "26.08.202510"),
    @(45895, 11, 3.441, "26.08.202511"),
    @(45895, 12, 3.732, "26.08.202512"),
    @(45895, 13, 3.726, "26.08.202513"),
    @(45895, 14, 3.513, "26.08.202514"),
    @(45895, 15, 3.738, "26.08.202515"),
    @(45895, 16, 3.831, "26.08.202516"),
    @(45895, 17, 3.398, "26.08.202517"),
    @(45895, 18, 2.474, "26.08.202518"),
    @(45895, 19, 1.325, "26.08.202519"),
    @(45895, 20, 0.364, "26.08.202520"),
    @(45895, 21, 0.027, "26.08.202521"),
    @(45895, 22, 0.054, "26.08.202522"),
    @(45895, 23, 0.06900000000000001, "26.08.202523"),
    @(45895, 24, 0.056, "26.08.202524"),
    @(45896, 1, 0.083, "27.08.20251"),
    @(45896, 2, 0.07199999999999999, "27.08.20252"),
    @(45896, 3, 0.114, "27.08.20253"),
    @(45896, 4, 0.108, "27.08.20254"),
    @(45896, 5, 0.104, "27.08.20255"),
    @(45896, 6, 0.099, "27.08.20256"),
    @(45896, 7, 0.06900000000000001, "27.08.20257"),
    @(45896, 8, 0.354, "27.08.20258"),
    @(45896, 9, 1.069, "27.08.20259"),
    @(45896, 10, 2.012, "27.08.202510"),
    @(45896, 11, 2.83, "27.08.202511"),
    @(45896, 12, 3.266, "27.08.202512"),
    @(45896, 13, 3.738, "27.08.202513"),
    @(45896, 14, 3.717, "27.08.202514"),
    @(45896, 15, 3.675, "27.08.202515"),
    @(45896, 16, 3.579, "27.08.202516"),
    @(45896, 17, 2.71, "27.08.202517"),
    @(45896, 18, 1.843, "27.08.202518"),
    @(45896, 19, 0.908, "27.08.202519"),
    @(45896, 20, 0.256, "27.08.202520"),
    @(45896, 21, 0.067, "27.08.202521"),
    @(45896, 22, 0.05, "27.08.202522"),
    @(45896, 23, 0.052, "27.08.202523"),
    @(45896, 24, 0.063, "27.08.202524"),
    @(45897, 1, 0.1, "28.08.20251"),
    @(45897, 2, 0.124, "28.08.20252"),
    @(45897, 3, 0.115, "28.08.20253"),
    @(45897, 4, 0.097, "28.08.20254"),
    @(45897, 5, 0.12, "28.08.20255"),
    @(45897, 6, 0.14, "28.08.20256"),
    @(45897, 7, 0.136, "28.08.20257"),
    @(45897, 8, 0.368, "28.08.20258"),
    @(45897, 9, 1.429, "28.08.20259"),
    @(45897, 10, 2.413, "28.08.202510"),
    @(45897, 11, 3.276, "28.08.202511"),
    @(45897, 12, 3.453, "28.08.202512"),
    @(45897, 13, 3.521, "28.08.202513"),
    @(45897, 14, 3.726, "28.08.202514"),
    @(45897, 15, 3.745, "28.08.202515"),
    @(45897, 16, 3.568, "28.08.202516"),
    @(45897, 17, 3.393, "28.08.202517")
)

$startRow = 2
for ($i = 0; $i -lt $forecast.Count; $i++) {
    $row = $startRow + $i
    $record = $forecast[$i]
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
}

